# TaskMaster/newreport.xlsx — "ToString и excel-report у BranchedTask"
#
# 1. Simple Tasks: rename the sample topic "myTopic" -> "chech"
# 2. Simple Tasks: append a new sample row (a sub-task, "InProcess")
# 3. Branched Tasks: populate the report with a "SubTasks" header column
#    plus the Simple Tasks headers, and one data row showing a task with
#    its sub-task underneath it.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Simple Tasks")
$ws2 = $wb.Worksheets.Item("Branched Tasks")

# --- Simple Tasks: rename myTopic -> chech -------------------------------
$ws1.Range("A2").Value = "chech"

# --- Simple Tasks: new row 3 (a sub task entry) --------------------------
$ws1.Range("A3").Value = "nya"
$ws1.Range("B3").Value = "g"
$ws1.Range("C3").Value = "InProcess"
$ws1.Range("G3").Value = "Valera"
$ws1.Range("H3").Value = "Valera"

# --- Branched Tasks: header row -------------------------------------------
$ws2.Range("A1").Value = "SubTasks"
$ws2.Range("B1").Value = "Topic"
$ws2.Range("C1").Value = "Description"
$ws2.Range("D1").Value = "State"
$ws2.Range("E1").Value = "Start"
$ws2.Range("F1").Value = "Finish"
$ws2.Range("G1").Value = "DeadLine"
$ws2.Range("H1").Value = "Performer"
$ws2.Range("I1").Value = "Owner"

# --- Branched Tasks: data row (task + nested sub task) --------------------
$ws2.Range("A2").Value = "chech"
$ws2.Range("B2").Value = "nya"
$ws2.Range("C2").Value = "g"
$ws2.Range("D2").Value = "InProcess"
$ws2.Range("H2").Value = "Valera"
$ws2.Range("I2").Value = "Valera"
